$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(213, 1).Value = "Soley Eczanesi"
$ws.Cells.Item(213, 2).Value = "`n+90 506 598 90 50"
$ws.Cells.Item(213, 3).Value = "`nBostancı, Kozyatağı, Seda Sk. No:11/A, 34742 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(213, 4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(213, 5).Value = "2026-02-11 01:47"

$ws.Cells.Item(214, 1).Value = "CADDE SAĞLIK ECZANESİ"
$ws.Cells.Item(214, 2).Value = "`n+90 216 356 00 85"
$ws.Cells.Item(214, 3).Value = "`nCADDE SAĞLIK ECZANESİ, Caddebostan, Bağdat Cad. NO:275 D:1, 34728 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(214, 4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(214, 5).Value = "2026-02-11 01:47"

$ws.Cells.Item(215, 1).Value = "Alp Eczanesi"
$ws.Cells.Item(215, 2).Value = "`n+90 546 128 82 46"
$ws.Cells.Item(215, 3).Value = "`nErenköy, Kamiller Sokağı No:5 D:B, 34738 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(215, 4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(215, 5).Value = "2026-02-11 01:47"

$ws.Cells.Item(216, 1).Value = "Karaca Eczanesi"
$ws.Cells.Item(216, 2).Value = "`n+90 216 748 08 88"
$ws.Cells.Item(216, 3).Value = "`n19 Mayıs, Oral Sk. No: 1/B, 34736 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(216, 4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(216, 5).Value = "2026-02-11 01:47"

$ws.Cells.Item(217, 1).Value = "GÖKÇEN ECZANESİ"
$ws.Cells.Item(217, 2).Value = "`n+90 216 629 86 06"
$ws.Cells.Item(217, 3).Value = "`nDumlupınar, Yazıcılar Sk., 34720 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(217, 4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(217, 5).Value = "2026-02-11 01:47"

$ws.Cells.Item(218, 1).Value = "Eczane Ertem"
$ws.Cells.Item(218, 2).Value = "`n+90 216 338 84 98"
$ws.Cells.Item(218, 3).Value = "`nCaferağa, Moda Cd. No:112, 34710 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(218, 4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(218, 5).Value = "2026-02-11 01:47"

$ws.Cells.Item(219, 1).Value = "Eczane Nihal"
$ws.Cells.Item(219, 2).Value = "`n+90 216 368 42 14"
$ws.Cells.Item(219, 3).Value = "`n19 Mayıs, Sinan Ercan Cd. No:30 C, 34736 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(219, 4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(219, 5).Value = "2026-02-11 01:47"

$ws.Cells.Item(220, 1).Value = "Eczane Saye"
$ws.Cells.Item(220, 2).Value = "`n+90 216 360 85 93"
$ws.Cells.Item(220, 3).Value = "`n19 Mayıs, Yıldız Sk. No:15/B, 34736 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(220, 4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(220, 5).Value = "2026-02-11 01:47"

$ws.Cells.Item(221, 1).Value = "Güleç Eczanesi"
$ws.Cells.Item(221, 2).Value = "`n+90 216 909 26 96"
$ws.Cells.Item(221, 3).Value = "`n19 Mayıs, Yıldız Sk. No:18/B, 34736 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(221, 4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(221, 5).Value = "2026-02-11 01:47"

$ws.Cells.Item(222, 1).Value = "Ergenekon Eczanesi"
$ws.Cells.Item(222, 2).Value = "`n+90 216 414 65 83"
$ws.Cells.Item(222, 3).Value = "`nZühtüpaşa, Kördere Sok. 22/A, 34724 Kadıköy/İstanbul, Türkiye"
$ws.Cells.Item(222, 4).Value = "Kadikoy/Istanbul"
$ws.Cells.Item(222, 5).Value = "2026-02-11 01:47"
